$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM line: microcontroller chip, plain text (no hyperlink) in C5
$ws.Range("C5").Value = "ATMega328P-PU"

# New BOM line: 9V battery clips, styled + hyperlinked like the other parts in C12
$ws.Range("C12").Value = "9V Battery Clips w/ Bare Wires"
$ws.Hyperlinks.Add($ws.Range("C12"), "https://www.sparkfun.com/products/9518")
$ws.Range("C12").Value = "9V Battery Clips w/ Bare Wires"
$ws.Range("C12").Style = "Hyperlink"

Write-Output "done"
